$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 58643.807797054498
$ws.Range("C2").Value = 32765.234998247099
$ws.Range("D2").Value = 1715.1757077961699

$ws.Range("B3").Value = 73395.0073361994
$ws.Range("C3").Value = 46523.187283071398
$ws.Range("D3").Value = 721.92845347583295

$ws.Range("B4").Value = 85104.733203326701
$ws.Range("C4").Value = 57788.213950614001
$ws.Range("D4").Value = 277.22925389091398

$ws.Range("B5").Value = 93821.9518878526
$ws.Range("C5").Value = 66342.553604917004
$ws.Range("D5").Value = 114.35022366793299

$ws.Range("B6").Value = 104210.498158174
$ws.Range("C6").Value = 76654.343338613995
$ws.Range("D6").Value = 37.593687041302999

$ws.Range("E19").Select()
